# Refresh the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns on
# Sheet1 with the latest scraped snapshot values. Cells keep their original
# text representation (e.g. "569.58", "  -0.55%  ") rather than being
# reinterpreted as numbers, matching how the source data is stored.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.530.23"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.491.27"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'569.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "'164.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.511"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").Value = "2.489.71"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "'0.352"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "2.945.76"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").Value = "69.394.11"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "'24.16"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("D18").Value = "2.494.65"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "'11.16"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("D20").Value = "'7.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").Value = "'347.34"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'3.88"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("E23").Value = "  -4.34%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'69.26"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("D27").Value = "2.617.49"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "'8.58"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.25%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "0.0₃0870"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").Value = "'7.57"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("E32").Value = "  -3.92%  "
$ws.Range("D33").Value = "'435.10"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.02%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").Value = "'155.53"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("E37").Value = "  -3.63%  "
$ws.Range("D38").Value = "'19.07"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "'18.11"
$ws.Range("D39").ClearFormats()
$ws.Range("D41").Value = "'0.313"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").Value = "'4.58"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.49%  "
$ws.Range("D43").Value = "'2.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +51.13%  "
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("E46").Value = "  -5.76%  "
$ws.Range("D47").Value = "'138.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").Value = "'3.42"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("E49").Value = "  -4.13%  "
$ws.Range("D50").Value = "'0.0723"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("E51").Value = "  -1.23%  "
